$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'89.261.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "'3.091.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'212.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("D6").Value = "'621.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("E7").Value = "  -6.37%  "
$ws.Range("D8").Value = "'0.828"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +16.57%  "
$ws.Range("D10").Value = "'3.089.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.37%  "
$ws.Range("E11").Value = "  +8.58%  "
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.01%  "
$ws.Range("D14").Value = "'5.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "'88.784.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").Value = "'32.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "'3.659.80"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").Value = "'3.100.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.29%  "
$ws.Range("D19").Value = "'3.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.94%  "
$ws.Range("D20").Value = "'0.0000212"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").Value = "'13.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.38%  "
$ws.Range("D22").Value = "'425.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("D23").Value = "'8.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("D24").Value = "'4.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D25").Value = "'5.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.71%  "
$ws.Range("D26").Value = "'12.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.86%  "
$ws.Range("D27").Value = "'83.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.50%  "
$ws.Range("D28").Value = "'3.255.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("E29").Value = "  +0.13%  "
$ws.Range("E30").Value = "  +6.89%  "
$ws.Range("D31").Value = "'0.166"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.84%  "
$ws.Range("D32").Value = "'8.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("D33").Value = "'507.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").Value = "'3.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.03%  "
$ws.Range("E35").Value = "  -3.13%  "
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("E37").Value = "  -4.22%  "
$ws.Range("D38").Value = "'22.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("D39").Value = "'22.27"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'0.364"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.137"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.81%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.62%  "
$ws.Range("D46").Value = "'145.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "'0.0697"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +14.82%  "
$ws.Range("D48").Value = "'43.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.30%  "
$ws.Range("E49").Value = "  +1.48%  "
$ws.Range("D50").Value = "'159.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.62%  "
$ws.Range("D51").Value = "'0.706"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.45%  "
